# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G ("K" - strikeouts) is re-derived from the updated source data
# (K counts instead of the previous Strike# metric) and the values are
# rewritten for each data row (rows 2-70).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(1,0,1,1,2,2,1,0,1,1,0,2,1,2,0,2,0,0,1,2,1,1,0,1,0,0,0,3,0,0,2,1,1,0,2,1,1,1,2,1,1,0,3,0,0,1,1,0,2,1,0,0,2,1,1,1,0,0,0,0,0,2,1,0,1,2,3,2,1)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = 2 + $i
    $ws.Range("G$row").Value = $newK[$i]
}
